$wb = $excel.ActiveWorkbook

# --- 1) Update "Estadisticos 2P" statistics (row 2 and row 5) ---
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")

$ws2.Cells.Item(2,4).Value = 28
$ws2.Cells.Item(2,5).Value = 28
$ws2.Cells.Item(2,6).Value = 8
$ws2.Cells.Item(2,7).Value = 22.22
$ws2.Cells.Item(2,8).Value = 8.800000000000001

$ws2.Cells.Item(5,4).Value = 10
$ws2.Cells.Item(5,5).Value = 10
$ws2.Cells.Item(5,6).Value = 12
$ws2.Cells.Item(5,7).Value = 54.55
$ws2.Cells.Item(5,8).Value = 8.6

# --- 2) Append new student rows to "Rescatables" sheet ---
$ws4 = $wb.Worksheets.Item("Rescatables")

# Column A - student id numbers
$ws4.Cells.Item(2,1).Value = 20330051920359
$ws4.Cells.Item(3,1).Value = 20330051920091
$ws4.Cells.Item(4,1).Value = 19330051920060
$ws4.Cells.Item(5,1).Value = 18330051920069
$ws4.Cells.Item(6,1).Value = 18330051920172
$ws4.Cells.Item(7,1).Value = 18330051920044

# Column B - Paterno (written fully before moving to next column so the
# shared-string table fills in the same order the source file used)
$ws4.Cells.Item(2,2).Value = "CASTILLO"
$ws4.Cells.Item(3,2).Value = "OLMOS"
$ws4.Cells.Item(4,2).Value = "HERNANDEZ"
$ws4.Cells.Item(5,2).Value = "MARIA"
$ws4.Cells.Item(6,2).Value = "RIOS"
$ws4.Cells.Item(7,2).Value = "ANTONIO"

# Column C - Materno
$ws4.Cells.Item(2,3).Value = "ROMERO"
$ws4.Cells.Item(3,3).Value = "CASTRO"
$ws4.Cells.Item(4,3).Value = "CID"
$ws4.Cells.Item(5,3).Value = "HERNANDEZ"
$ws4.Cells.Item(6,3).Value = "OCHOA"
$ws4.Cells.Item(7,3).Value = "GAMBINO"

# Column D - Nombres
$ws4.Cells.Item(2,4).Value = "KARLA JOVANA"
$ws4.Cells.Item(3,4).Value = "ANGEL"
$ws4.Cells.Item(4,4).Value = "ELIUTH ADELFO"
$ws4.Cells.Item(5,4).Value = "AMALIO JARET"
$ws4.Cells.Item(6,4).Value = "JONATHAN"
$ws4.Cells.Item(7,4).Value = "DAVID"

# Column E - Nombre_Largo (module name)
$ws4.Cells.Item(2,5).Value = "REALIZA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(3,5).Value = "REALIZA INSTALACIONES ELÉCTRICAS"
$ws4.Cells.Item(4,5).Value = "MANTIENE EN OPERACIÓN CIRCUITOS DE CONTROL ELECTRÓNICO"
$ws4.Cells.Item(5,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(6,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"
$ws4.Cells.Item(7,5).Value = "REALIZA MANTENIMIENTO EN EL SISTEMA DE DISTRIBUCIÓN DE ENERGÍA ELÉCTRICA"

# Column F - Grupo
$ws4.Cells.Item(2,6).Value = "2AEV"
$ws4.Cells.Item(3,6).Value = "2AEV"
$ws4.Cells.Item(4,6).Value = "4AEV"
$ws4.Cells.Item(5,6).Value = "6AEM"
$ws4.Cells.Item(6,6).Value = "6AEM"
$ws4.Cells.Item(7,6).Value = "6BEV"

# Column G - Reprobadas
$ws4.Cells.Item(2,7).Value = 2
$ws4.Cells.Item(3,7).Value = 2
$ws4.Cells.Item(4,7).Value = 2
$ws4.Cells.Item(5,7).Value = 2
$ws4.Cells.Item(6,7).Value = 2
$ws4.Cells.Item(7,7).Value = 2
